$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unnecessary "Got" column (column C), which duplicated the
# "Expected" column (B) for every test case row.
$ws.Range("C1:C11").ClearContents()
